$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): F5 11457 -> 11460, F9 11403 -> 11405
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 11460
$wsExhibit.Range("F9").Value = 11405

# Sheet "全部类型" (sheet4): F7 11457 -> 11460, F11 11403 -> 11405
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 11460
$wsAll.Range("F11").Value = 11405
